$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.143.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "'3.058.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'389.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "'101.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").Value = "'0.531"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.30%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.581"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("D10").Value = "'36.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "'0.0845"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").Value = "'3.537.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "'18.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").Value = "'7.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "'3.054.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "'1.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("D18").Value = "'10.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "'51.163.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "'3.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("D21").Value = "'12.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "'0.0₃0953"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").Value = "'69.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "'263.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "'3.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").Value = "'7.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.35%  "
$ws.Range("D27").Value = "'26.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'7.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.63%  "
$ws.Range("D30").Value = "'0.163"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.81%  "
$ws.Range("D31").Value = "'0.105"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").Value = "'10.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.00%  "
$ws.Range("D33").Value = "'0.0487"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.69%  "
$ws.Range("D34").Value = "'35.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.66%  "
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").Value = "'49.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "'3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "'0.290"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'128.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").Value = "'16.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("D42").Value = "'1.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("D43").Value = "'0.115"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'2.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'3.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("D46").Value = "'21.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("D48").Value = "'2.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("D49").Value = "'2.063.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("B50").Value = "BEAM"
$ws.Range("C50").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D50").Value = "'0.0316"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("B51").Value = "FlareNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/2hOSU_JYX+flarenetwork-flr"
$ws.Range("D51").Value = "'0.0471"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +18.84%  "
